# Auto-update stock values: 2025-12-09 07:55:33 UTC
#
# Appends a new trading-day column (2025-12-08) to each of the 12 data
# sheets (시가/고가/저가/종가/거래량/s20/s60/z20/z60/gap/std/quant).
# Each sheet gets one new column, one column to the right of its current
# last column, containing:
#   row 1 -> date 20251208 (header, same style as the previous header cell)
#   row 2 -> new value for ticker 1 (INVESCO QQQ TRUST)
#   row 3 -> new value for ticker 2 (PROSHARES QQQ 3X)
#
# Some sheets store the header date as a genuine number (t="n"); others
# (gap / std / quant) store it as text (t="inlineStr" in the original,
# re-serialised as a shared string by this engine) even though the cell
# format is plain "General". To reproduce a *text* cell without altering
# the cell's number format (and therefore without minting a brand new
# style index), we briefly stash a text formula in the cell and then
# freeze it back down to a literal value with Paste Special -> Values.

function Add-DayColumn {
    param(
        $ws,
        [int]$SrcCol,
        [int]$DstCol,
        [string]$HeaderValue,
        [bool]$HeaderIsText,
        $Row2Value,
        $Row3Value,
        [double]$RawColWidth
    )

    # 1) Clone the format of the previous header column's cells onto the
    #    whole new column (via the header cell) so font/fill match exactly.
    $ws.Cells.Item(1, $SrcCol).Copy()
    $ws.Cells.Item(1, $DstCol).PasteSpecial(-4122)  # xlPasteFormats

    # 2) Header cell (row 1): date value, numeric or text depending on sheet.
    if ($HeaderIsText) {
        $ws.Cells.Item(1, $DstCol).Formula = "=""$HeaderValue"""
        $ws.Cells.Item(1, $DstCol).Copy()
        $ws.Cells.Item(1, $DstCol).PasteSpecial(-4163)  # xlPasteValues
    } else {
        $ws.Cells.Item(1, $DstCol).Value = [double]$HeaderValue
    }

    # 3) Data rows.
    $ws.Cells.Item(2, $DstCol).Value = $Row2Value
    $ws.Cells.Item(3, $DstCol).Value = $Row3Value

    # 4) Column width matching the sheet's existing data columns.
    $ws.Columns.Item($DstCol).ColumnWidth = $RawColWidth
}

$wb = $excel.ActiveWorkbook

$narrow = 9.166666666666666   # raw <col width="10">
$wide   = 11.166666666666666  # raw <col width="12">

# 시가 (open) - col 71 (BS) -> 72 (BT)
Add-DayColumn $wb.Worksheets.Item("시가") 71 72 "20251208" $false 627.21 56.57 $wide

# 고가 (high) - col 71 (BS) -> 72 (BT)
Add-DayColumn $wb.Worksheets.Item("고가") 71 72 "20251208" $false 628.84 57.04 $wide

# 저가 (low) - col 71 (BS) -> 72 (BT)
Add-DayColumn $wb.Worksheets.Item("저가") 71 72 "20251208" $false 621.69 55.12 $wide

# 종가 (close) - col 71 (BS) -> 72 (BT)
Add-DayColumn $wb.Worksheets.Item("종가") 71 72 "20251208" $false 624.28 55.8 $wide

# 거래량 (volume) - col 71 (BS) -> 72 (BT)
Add-DayColumn $wb.Worksheets.Item("거래량") 71 72 "20251208" $false 43462406 69700018 $wide

# s20 - col 52 (AZ) -> 53 (BA)
Add-DayColumn $wb.Worksheets.Item("s20") 52 53 "20251208" $false 97 14 $narrow

# s60 - col 12 (L) -> 13 (M)
Add-DayColumn $wb.Worksheets.Item("s60") 12 13 "20251208" $false 77 13 $narrow

# z20 - col 52 (AZ) -> 53 (BA)
Add-DayColumn $wb.Worksheets.Item("z20") 52 53 "20251208" $false 50 -34 $narrow

# z60 - col 12 (L) -> 13 (M)
Add-DayColumn $wb.Worksheets.Item("z60") 12 13 "20251208" $false 61 -90 $narrow

# gap - col 52 (AZ) -> 53 (BA) ; header stored as text
Add-DayColumn $wb.Worksheets.Item("gap") 52 53 "20251208" $true 102 75 $wide

# std - col 33 (AG) -> 34 (AH) ; header stored as text
Add-DayColumn $wb.Worksheets.Item("std") 33 34 "20251208" $true 2.99 52.67 $wide

# quant - col 12 (L) -> 13 (M) ; header stored as text
Add-DayColumn $wb.Worksheets.Item("quant") 12 13 "20251208" $true 37 52 $wide
